$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,(-0.2281722784908007,-0.9450025826423425)
    ,(-0.4327561053543218,-1.007418635742523)
    ,(-0.4807351009405674,-1.487200379883007)
    ,(-0.1698080811235126,-0.8112142702687729)
    ,(-0.1435722966780263,0.6600759463522082)
    ,(-0.3087582281763757,-0.9461878793287239)
    ,(-0.4587794576881659,-0.7684964777525181)
    ,(-0.5757145459952729,-1.462759555691188)
    ,(-0.4514694563417559,-0.5340374614123443)
    ,(-0.2942052311451911,-1.144042715997361)
    ,(-0.2247271485605528,-0.7656606421897426)
    ,(-0.06860114684896625,-0.1882775918802016)
    ,(-0.2538971980255557,-0.2197522014503121)
    ,(-0.1166161989110121,0.0815373402541214)
    ,(-0.1701404519061747,-0.1923993313668516)
    ,(0.08573270739481063,0.5588443610890836)
    ,(0.0193536880717859,0.392549337515866)
    ,(0.02971326152010036,0.4603559120995097)
    ,(-0.1192912577763853,0.0547005355363147)
    ,(0.03902512980374093,0.3346442390706392)
    ,(0.06278842316261118,0.7322812204539821)
    ,(0.03385708262665894,0.3933558330928552)
    ,(0.8112252719489321,1.627973688413401)
    ,(0.1302275998330896,0.7456943197896312)
    ,(0.1515711138524058,0.7094818840958533)
    ,(0.105521761749875,0.7276387023418044)
    ,(0.2982722637802794,0.9807189957634788)
    ,(0.6314185549165483,1.690780017481349)
    ,(0.2013955421666478,0.7860223835656993)
    ,(0.02483719257675755,0.5980581981659435)
    ,(0.1736746342570804,0.920810253958359)
    ,(0.1182409424342906,0.810598167429427)
    ,(0.08883829216080727,0.4330117839264488)
    ,(0.4017800659878021,1.328501186268377)
    ,(0.2607131313157371,0.9701650371360311)
    ,(0.06583600360424673,0.3123393246518008)
    ,(0.3322697090984311,1.704254102429345)
    ,(-0.07060304971567119,-0.1644465890263914)
    ,(0.1786626385972376,0.9409548739025546)
    ,(-0.1167628966199462,0.65561474377252)
    ,(0.2418059993966016,1.25461173782584)
    ,(0.1916608245638603,1.117694387931033)
    ,(-0.1214147904016063,-0.1132243179000224)
    ,(-0.1057977733444171,0.1640665209056412)
    ,(-0.1894767094903146,-0.5233256499321619)
    ,(-0.1889778744099086,-0.4923663570447318)
    ,(-0.2285539900582748,-0.5456421001190956)
    ,(-0.2211376028321234,-0.5934885442690832)
    ,(-0.1542317080930548,-0.3837213988737411)
    ,(-0.2132182103331804,-0.6540461911515563)
    ,(-0.2132182103331804,-0.6540461911515563)
    ,(-0.2071250119793882,-0.4530253688300155)
    ,(-0.207832779078593,-0.549527211844548)
    ,(-0.1710718357200046,-0.3492374008156431)
    ,(-0.1500580597870149,-0.4337032655013847)
    ,(-0.1976938154196263,-0.5293200855273208)
    ,(-0.1803430406994986,-0.588202982933294)
    ,(-0.2214572396297638,-0.6682422967837788)
    ,(-0.2093823238006601,-0.6393892857454495)
    ,(-0.2654826144233455,-0.3862114602999622)
    ,(-0.1642097265672183,-0.07146077379300829)
    ,(-0.3436445555912591,-1.135682047689604)
    ,(-0.2722895532434333,-0.6684407944034956)
    ,(-0.193627150722176,-0.6654233704590338)
    ,(-0.1196660269786407,-0.12790649042938)
    ,(-0.05897642808261178,-0.2848051350889697)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}